$d = $word.ActiveDocument

# wdFindContinue = 1, wdReplaceOne = 1, wdReplaceAll = 2 (used below)

# 1st person: "1ª persona ():" -> "1ª persona (Augusto):"
$d.Content.Find.Execute("persona ():", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "persona (Augusto):", 2)

# 2nd person: "2ª persona():" -> "2ª persona(Jose):"
$d.Content.Find.Execute("2ª persona():", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "2ª persona(Jose):", 2)

# 3rd person: "3ª persona(    ):" (four spaces) -> "3ª persona( Tomas):"
$d.Content.Find.Execute("3ª persona(    ):", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "3ª persona( Tomas):", 2)

# 4th person: "4ª persona():" -> "4ª persona(Alberto):"
$d.Content.Find.Execute("4ª persona():", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "4ª persona(Alberto):", 2)

# 5th person: "5ª persona():" -> "5ª persona(Adrian):"
$d.Content.Find.Execute("5ª persona():", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "5ª persona(Adrian):", 2)

# The "_GoBack" bookmark used to sit right before "3.1 y 3.2" (i.e. right after
# "Puntos: " in the 4th person's points paragraph). In the edited document it
# instead sits right after the newly-inserted "Alberto" name (immediately
# before the closing "):"). Re-anchor it there.
$rng = $d.Content
$rng.Find.Execute("Alberto")
$rng.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rng)
